# This change updates the NATMI Tgm2-Sdc4 ligand-receptor pair results
# (columns E..T for data rows 2..17) to reflect a recomputation where the
# number of ligand-/receptor-expressing cells (columns E and K) increased
# from 1 to 3, changing all of the downstream expression / specificity /
# edge-weight statistics. Columns A-D (Sending cluster, Ligand symbol,
# Receptor symbol, Target cluster) and F, L (detection rates) are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array holds the new values, in column order, for:
# E (Ligand-expressing cells), F (Ligand detection rate),
# G (Ligand average expr), H (Ligand total expr),
# I (Ligand specificity avg), J (Ligand specificity total),
# K (Receptor-expressing cells), L (Receptor detection rate),
# M (Receptor average expr), N (Receptor total expr),
# O (Receptor specificity avg), P (Receptor specificity total),
# Q (Edge average expr weight), R (Edge total expr weight),
# S (Edge avg expr specificity), T (Edge total expr specificity)
$rows = @(
    @(3,1,79.612463,238.837389,0.6728436998494041,0.6728436998494042,3,1,3.001642333333333,9.004927,0.05169795991651582,0.05169795991651582,238.9681392017337,2150.713252815603,0.03478464662489469,0.0347846466248947),
    @(3,1,79.612463,238.837389,0.6728436998494041,0.6728436998494042,3,1,11.05428166666667,33.162845,0.1903903753498087,0.1903903753498088,880.0585901790785,7920.527311611706,0.1281029645660821,0.1281029645660821),
    @(3,1,79.612463,238.837389,0.6728436998494041,0.6728436998494042,3,1,16.49405866666666,49.482176,0.2840808761059341,0.2840808761059341,1313.132635319829,11818.19371787846,0.1911420277355768,0.1911420277355769),
    @(3,1,79.612463,238.837389,0.6728436998494041,0.6728436998494042,3,1,27.51115433333333,82.533463,0.4738307886277414,0.4738307886277414,2190.23075644979,19712.07680804811,0.3188140609228504,0.3188140609228505),
    @(3,1,3.815058666666667,11.445176,0.03224291890608301,0.03224291890608302,3,1,3.001642333333333,9.004927,0.05169795991651582,0.05169795991651582,11.45144159801689,103.062974382152,0.001666893129198149,0.00166689312919815),
    @(3,1,3.815058666666667,11.445176,0.03224291890608301,0.03224291890608302,3,1,11.05428166666667,33.162845,0.1903903753498087,0.1903903753498088,42.17273307619112,379.5545976857201,0.006138741432902588,0.006138741432902591),
    @(3,1,3.815058666666667,11.445176,0.03224291890608301,0.03224291890608302,3,1,16.49405866666666,49.482176,0.2840808761059341,0.2840808761059341,62.92580146477511,566.3322131829759,0.009159596651052647,0.00915959665105265),
    @(3,1,3.815058666666667,11.445176,0.03224291890608301,0.03224291890608302,3,1,27.51115433333333,82.533463,0.4738307886277414,0.4738307886277414,104.9566677693876,944.6100099244879,0.01527768769292963,0.01527768769292963),
    @(3,1,30.45313966666667,91.359419,0.2573743154429307,0.2573743154429307,3,1,3.001642333333333,9.004927,0.05169795991651582,0.05169795991651582,91.40943320637923,822.684898857413,0.01330572704330933,0.01330572704330933),
    @(3,1,30.45313966666667,91.359419,0.2573743154429307,0.2573743154429307,3,1,11.05428166666667,33.162845,0.1903903753498087,0.1903903753498088,336.6375835096728,3029.738251587055,0.04900159252257964,0.04900159252257966),
    @(3,1,30.45313966666667,91.359419,0.2573743154429307,0.2573743154429307,3,1,16.49405866666666,49.482176,0.2840808761059341,0.2840808761059341,502.2958722461937,4520.662850215744,0.07311512101819279,0.07311512101819281),
    @(3,1,30.45313966666667,91.359419,0.2573743154429307,0.2573743154429307,3,1,27.51115433333333,82.533463,0.4738307886277414,0.4738307886277414,837.801025304222,7540.209227737997,0.1219518748588489,0.121951874858849),
    @(3,1,4.441711333333334,13.325134,0.03753906580158222,0.03753906580158223,3,1,3.001642333333333,9.004927,0.05169795991651582,0.05169795991651582,13.33242877057978,119.991858935218,0.001940693119113647,0.001940693119113648),
    @(3,1,4.441711333333334,13.325134,0.03753906580158222,0.03753906580158223,3,1,11.05428166666667,33.162845,0.1903903753498087,0.1903903753498088,49.09992816069224,441.8993534462301,0.007147076828244408,0.00714707682824441),
    @(3,1,4.441711333333334,13.325134,0.03753906580158222,0.03753906580158223,3,1,16.49405866666666,49.482176,0.2840808761059341,0.2840808761059341,73.26184731239822,659.3566258115841,0.01066413070111179,0.01066413070111179),
    @(3,1,4.441711333333334,13.325134,0.03753906580158222,0.03753906580158223,3,1,27.51115433333333,82.533463,0.4738307886277414,0.4738307886277414,122.1966059954491,1099.769453959042,0.01778716515311238,0.01778716515311238)
)

$numRows = $rows.Count
$numCols = $rows[0].Count

$arr = New-Object 'object[,]' $numRows,$numCols
for ($i = 0; $i -lt $numRows; $i++) {
    for ($j = 0; $j -lt $numCols; $j++) {
        $arr[$i,$j] = $rows[$i][$j]
    }
}

$ws.Range("E2:T17").Value = $arr
